{"js": "// Update the 100 arithmetic-answer cells in the document's table (20 rows x 5\n// columns) to the new set of problems, applied positionally in document\n// (row-major) order so it matches the original cells 1:1 even where some old\n// values repeat (e.g. \"6+35=41\" appears more than once but maps to different\n// new values depending on its position).\nconst newValues = [\n  [\"62-38=24\", \"29+35=64\", \"8+66=74\", \"61-12=49\", \"28+24=52\"],\n  [\"67-39=28\", \"82-54=28\", \"56-8=48\", \"28+46=74\", \"68+9=77\"],\n  [\"28+33=61\", \"23+9=32\", \"38+17=55\", \"17+68=85\", \"8+5=13\"],\n  [\"4+69=73\", \"95-9=86\", \"67+28=95\", \"73-58=15\", \"29+9=38\"],\n  [\"55-8=47\", \"28+5=33\", \"48+15=63\", \"36-18=18\", \"76-9=67\"],\n  [\"74-7=67\", \"64-57=7\", \"91-2=89\", \"73-34=39\", \"72-43=29\"],\n  [\"27+14=41\", \"84-56=28\", \"55+29=84\", \"12+69=81\", \"72-53=19\"],\n  [\"28+38=66\", \"37+14=51\", \"80-48=32\", \"15+27=42\", \"13+38=51\"],\n  [\"53+19=72\", \"91-23=68\", \"50-32=18\", \"35+58=93\", \"81-4=77\"],\n  [\"69+6=75\", \"88+7=95\", \"90-29=61\", \"81-65=16\", \"9+29=38\"],\n  [\"28+58=86\", \"26+36=62\", \"59+8=67\", \"13+8=21\", \"17+6=23\"],\n  [\"44-18=26\", \"66-47=19\", \"35+18=53\", \"9+13=22\", \"49+16=65\"],\n  [\"50-1=49\", \"26+67=93\", \"6+5=11\", \"21-18=3\", \"85-38=47\"],\n  [\"14-9=5\", \"64+18=82\", \"94-85=9\", \"90-27=63\", \"42-15=27\"],\n  [\"9+33=42\", \"31-6=25\", \"53-4=49\", \"18+39=57\", \"60-16=44\"],\n  [\"38+3=41\", \"29+38=67\", \"55-48=7\", \"93-64=29\", \"66-48=18\"],\n  [\"24-9=15\", \"18+74=92\", \"92-78=14\", \"25+38=63\", \"25+68=93\"],\n  [\"94-78=16\", \"65+29=94\", \"80-46=34\", \"73-66=7\", \"91-24=67\"],\n  [\"80-66=14\", \"32-15=17\", \"40-25=15\", \"49+49=98\", \"42-39=3\"],\n  [\"29+7=36\", \"91-9=82\", \"33+49=82\", \"18+33=51\", \"19+32=51\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst actualCols = table.values && table.values[0] ? table.values[0].length : 0;\n\nif (table.rowCount === newValues.length && actualCols === newValues[0].length) {\n  // Assign the full 2D array in one shot; this keeps each cell's existing\n  // paragraph/run formatting (font, size, alignment) intact and only swaps\n  // the text content, matching how the diff only touches the <w:t> nodes.\n  table.values = newValues;\n  await context.sync();\n} else {\n  // Fallback: table shape differs from what we expect, so update cell by\n  // cell (still positional/row-major) using a range replace that preserves\n  // each cell's existing formatting.\n  const cols = newValues[0].length;\n  for (let r = 0; r < newValues.length; r++) {\n    for (let c = 0; c < cols; c++) {\n      const cell = table.getCell(r, c);\n      const range = cell.body.getRange(\"Whole\");\n      range.insertText(newValues[r][c], Word.InsertLocation.replace);\n    }\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 100 arithmetic-answer cells in the table (20 rows x 5 columns)\n# to the new set of problems, applied positionally in document (row-major) order,\n# matching the row-major order of the original cells.\n$d = $word.ActiveDocument\n\n$newValues = @(\n    \"62-38=24\",\n    \"29+35=64\",\n    \"8+66=74\",\n    \"61-12=49\",\n    \"28+24=52\",\n    \"67-39=28\",\n    \"82-54=28\",\n    \"56-8=48\",\n    \"28+46=74\",\n    \"68+9=77\",\n    \"28+33=61\",\n    \"23+9=32\",\n    \"38+17=55\",\n    \"17+68=85\",\n    \"8+5=13\",\n    \"4+69=73\",\n    \"95-9=86\",\n    \"67+28=95\",\n    \"73-58=15\",\n    \"29+9=38\",\n    \"55-8=47\",\n    \"28+5=33\",\n    \"48+15=63\",\n    \"36-18=18\",\n    \"76-9=67\",\n    \"74-7=67\",\n    \"64-57=7\",\n    \"91-2=89\",\n    \"73-34=39\",\n    \"72-43=29\",\n    \"27+14=41\",\n    \"84-56=28\",\n    \"55+29=84\",\n    \"12+69=81\",\n    \"72-53=19\",\n    \"28+38=66\",\n    \"37+14=51\",\n    \"80-48=32\",\n    \"15+27=42\",\n    \"13+38=51\",\n    \"53+19=72\",\n    \"91-23=68\",\n    \"50-32=18\",\n    \"35+58=93\",\n    \"81-4=77\",\n    \"69+6=75\",\n    \"88+7=95\",\n    \"90-29=61\",\n    \"81-65=16\",\n    \"9+29=38\",\n    \"28+58=86\",\n    \"26+36=62\",\n    \"59+8=67\",\n    \"13+8=21\",\n    \"17+6=23\",\n    \"44-18=26\",\n    \"66-47=19\",\n    \"35+18=53\",\n    \"9+13=22\",\n    \"49+16=65\",\n    \"50-1=49\",\n    \"26+67=93\",\n    \"6+5=11\",\n    \"21-18=3\",\n    \"85-38=47\",\n    \"14-9=5\",\n    \"64+18=82\",\n    \"94-85=9\",\n    \"90-27=63\",\n    \"42-15=27\",\n    \"9+33=42\",\n    \"31-6=25\",\n    \"53-4=49\",\n    \"18+39=57\",\n    \"60-16=44\",\n    \"38+3=41\",\n    \"29+38=67\",\n    \"55-48=7\",\n    \"93-64=29\",\n    \"66-48=18\",\n    \"24-9=15\",\n    \"18+74=92\",\n    \"92-78=14\",\n    \"25+38=63\",\n    \"25+68=93\",\n    \"94-78=16\",\n    \"65+29=94\",\n    \"80-46=34\",\n    \"73-66=7\",\n    \"91-24=67\",\n    \"80-66=14\",\n    \"32-15=17\",\n    \"40-25=15\",\n    \"49+49=98\",\n    \"42-39=3\",\n    \"29+7=36\",\n    \"91-9=82\",\n    \"33+49=82\",\n    \"18+33=51\",\n    \"19+32=51\"\n)\n\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($i -ge $newValues.Count) { continue }\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i++\n    }\n}\n"}
